$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Start"
$ws.Range("B1").Value = "End"
$ws.Range("C1").Value = "Task Description"
$ws.Range("D1").Value = "Duration"

# --- Logged task rows -------------------------------------------------------
$ws.Range("A2").Value = "2021-08-26 13:11:14"
$ws.Range("B2").Value = "2021-08-26 13:11:14"
$ws.Range("C2").Value = "jakiś Task"
$ws.Range("D2").Value = "0:00:00"

$ws.Range("A3").Value = "2021-08-26 13:12:42"
$ws.Range("B3").Value = "2021-08-26 13:12:42"
$ws.Range("C3").Value = "sffasfsf"
$ws.Range("D3").Value = "0:00:00"

$ws.Range("A4").Value = "2021-08-26 13:13:36"
$ws.Range("B4").Value = "2021-08-26 13:13:36"
$ws.Range("C4").Value = "rgdfgfdsg"
$ws.Range("D4").Value = "0:00:00"

$ws.Range("A5").Value = "2021-08-26 13:15:48"
$ws.Range("B5").Value = "2021-08-26 13:15:48"
$ws.Range("C5").Value = "wdfgdfg"
$ws.Range("D5").Value = "0:00:00"

$ws.Range("A6").Value = "2021-08-26 13:21:47"
$ws.Range("B6").Value = "2021-08-26 13:22:10"
$ws.Range("C6").Value = "adgdfga"
$ws.Range("D6").Value = "0:00:23"

$ws.Range("A7").Value = "2021-08-26 13:25:31"
$ws.Range("B7").Value = "2021-08-26 13:25:33"
$ws.Range("C7").Value = "jakiś opis"
$ws.Range("D7").Value = "0:00:01"

$ws.Range("A8").Value = "2021-08-26 13:30:35"
$ws.Range("B8").Value = "2021-08-26 13:30:47"
$ws.Range("C8").Value = "przykładowy opis"
$ws.Range("D8").Value = "0:00:11"

$ws.Range("A9").Value = "2021-08-26 13:40:48"
$ws.Range("B9").Value = "2021-08-26 13:40:53"
$ws.Range("C9").Value = "costam"
$ws.Range("D9").Value = "0:00:04"

# --- Selection / view state -------------------------------------------------
$ws.Range("A1:D3").Select() | Out-Null
